# "added fib 4.c (gold)" - append a third algorithm variant (fib_gold_4c)
# results to the pytest-results log worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a plain default-styled TEXT value into a cell, even when
# the text looks numeric/dateish (so Excel doesn't silently convert it to
# a number/date). We force text via a leading apostrophe, then paste the
# *format only* from a clean, never-touched default cell ($ref) on top so
# the cell doesn't end up tagged with a stray "quote-prefixed" style.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text, $ref) {
    $cell.Value = "'" + $text
    $ref.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# A completely untouched, default-styled cell we can borrow clean
# formatting from throughout the script.
$blank = $ws.Range("Z500")

# ---------------------------------------------------------------------
# 1. Drop the old rows 23-41: this removed the duplicated "echo" block
#    (which just repeated fib_iteration_4b's detail rows as one long
#    comma joined string per row) together with the separate detail
#    table for fib_iteration_4b that used to live at rows 33-41.
# ---------------------------------------------------------------------
$ws.Rows("23:41").Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. Summary table (rows 2-9): add a new column D with fib_gold_4c
#    timings alongside the existing fib_recursion_4a / fib_iteration_4b
#    columns.
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "fib_gold_4c"

$goldTimes = @("0.0018188","0.0025526","0.0061877","0.0050125","0.001975","0.0049472","0.0033072")
for ($i = 0; $i -lt 7; $i++) {
    $row = 3 + $i
    Set-TextValue $ws.Cells.Item($row, 4) $goldTimes[$i] $blank
}

# ---------------------------------------------------------------------
# 3. Column widths: column A now fits the longer "fib_gold_4c,test..."
#    strings, column D fits the new gold timing column.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 27.1666666666667
$ws.Columns("D:D").ColumnWidth = 9.52994791666667

# ---------------------------------------------------------------------
# 4. Rebuild rows 23-40.
#
#    rows 23 / 31          : plain "----" divider lines (default style)
#    rows 24-30             : fib_iteration_4b detail rows (promoted out
#                              of the old echo block into real columns)
#    rows 33-39             : fib_gold_4c echo block (one combined
#                              string per test, "console log" style)
#    row 40                 : plain "----" divider line (echo style)
# ---------------------------------------------------------------------

$divider = "------------------------------------------------------------------------"

Set-TextValue $ws.Range("A23") $divider $blank

$files = @("test.0.in","test.1.in","test.2.in","test.3.in","test.4.in","test.5.in","test.6.in")
$iterTimes = @("0.003152","0.0078031","0.0032145","0.0032896","0.002602","0.005845","0.0030685")
$iterStamps = @(
    "2020-02-07 17:59:08.706976",
    "2020-02-07 17:59:08.714354",
    "2020-02-07 17:59:08.717425",
    "2020-02-07 17:59:08.721369",
    "2020-02-07 17:59:08.724383",
    "2020-02-07 17:59:08.730399",
    "2020-02-07 17:59:08.733405"
)

for ($i = 0; $i -lt 7; $i++) {
    $row = 24 + $i
    Set-TextValue $ws.Cells.Item($row, 1) "fib_iteration_4b" $blank
    Set-TextValue $ws.Cells.Item($row, 2) $files[$i] $blank
    $ws.Cells.Item($row, 3).Value = $i
    $ws.Cells.Item($row, 4).Value = $true
    Set-TextValue $ws.Cells.Item($row, 5) $iterTimes[$i] $blank
    Set-TextValue $ws.Cells.Item($row, 6) $iterStamps[$i] $blank
}

Set-TextValue $ws.Range("A31") $divider $blank

# Console-echo style block for the new fib_gold_4c run: reuse the same
# "Consolas 10pt, vertically centered" look previously used for this
# kind of divider/echo block.
$echoStyleSource = $ws.Range("A23")

$goldStamps = @(
    "2020-02-07 18:10:27.938622",
    "2020-02-07 18:10:27.941126",
    "2020-02-07 18:10:27.947649",
    "2020-02-07 18:10:27.952662",
    "2020-02-07 18:10:27.954668",
    "2020-02-07 18:10:27.959682",
    "2020-02-07 18:10:27.963061"
)

for ($i = 0; $i -lt 7; $i++) {
    $row = 33 + $i
    $line = "fib_gold_4c,test." + $i + ".in," + $i + ",True," + $goldTimes[$i] + "," + $goldStamps[$i]
    Set-TextValue $ws.Cells.Item($row, 1) $line $blank
    $ws.Cells.Item($row, 1).Font.Name = "Consolas"
    $ws.Cells.Item($row, 1).Font.Size = 10
    $ws.Cells.Item($row, 1).VerticalAlignment = -4108
}

Set-TextValue $ws.Range("A40") $divider $blank
$ws.Range("A40").Font.Name = "Consolas"
$ws.Range("A40").Font.Size = 10
$ws.Range("A40").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Clear the stray Z500 scratch cell we used as a clean formatting
#    source, and restore the selection the way the author left it.
# ---------------------------------------------------------------------
$blank.Clear() | Out-Null
$ws.Range("C36").Select() | Out-Null
